{"js": "// Correct ES header for i21-23\n// The Spanish paragraph text had a garbled/incorrect translation; replace it\n// with the corrected wording. The trailing colon is left as its own run,\n// matching how the paragraph is split in the corrected document.\nconst oldTextNoColon =\n  \"Si alguien dice algo que mi hijo/a piensa que est\u00e1 equivocado e es malo, mi hijo/s se siente asustado/a diciendo o que piensa si esa persona es\";\nconst newTextNoColon =\n  \"Si alguien dice algo que mi hijo/a piensa que es equivocado o malo, mi hijo/a se siente asustado/a de decir lo que piensa si esa persona es\";\n\nconst results = context.document.body.search(oldTextNoColon, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found\");\n}\n\nresults.items[0].insertText(newTextNoColon, \"Replace\");\nawait context.sync();\n", "ps1": "# Correct ES header for i21-23\n# The Spanish paragraph text had a garbled/incorrect translation; replace it\n# with the corrected wording.\n$d = $word.ActiveDocument\n\n$oldText = \"Si alguien dice algo que mi hijo/a piensa que est\u00e1 equivocado e es malo, mi hijo/s se siente asustado/a diciendo o que piensa si esa persona es:\"\n$newText = \"Si alguien dice algo que mi hijo/a piensa que es equivocado o malo, mi hijo/a se siente asustado/a de decir lo que piensa si esa persona es:\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
